# Highlight specific checklist paragraphs, matching the reference revision.
# Word stamps the paragraph-mark run (w:pPr/w:rPr) with the same highlight
# whenever a whole paragraph (incl. its pilcrow) is highlighted, so we rebuild
# each target paragraph's OOXML (preserving its existing runs/proofErr marks)
# with a <w:highlight> added to the pPr and to every run, then push it back in
# with Range.InsertXML (InsertXML replaces the exact range it is called on).

$d = $word.ActiveDocument

# Paragraph 43: Cubes are drawn using a custom WorldOpaque shader program... -> red
$paraXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="01A9BA26" w14:textId="77777777" w:rsidR="00D361C0" w:rsidRDefault="00D361C0" w:rsidP="007A28AD"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="red"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>Cubes are drawn using a custom “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>WorldOpaque</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="red"/></w:rPr><w:t>” shader program which does something per-pixel to demonstrate its use</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(43).Range.InsertXML($paraXml)

# Paragraph 72: (5 points) XML support -> cyan
$paraXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="398FE059" w14:textId="77777777" w:rsidR="0067117D" w:rsidRDefault="00616B75" w:rsidP="00616B75"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidR="00705E8E"><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t xml:space="preserve"> points) </w:t></w:r><w:r w:rsidR="00CF7E20"><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>XML</w:t></w:r><w:r w:rsidR="00CC1722"><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t xml:space="preserve"> support</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(72).Range.InsertXML($paraXml)

# Paragraph 73: You must still be able to open and parse data elements... -> cyan
$paraXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5F9E93B1" w14:textId="77777777" w:rsidR="00CC1722" w:rsidRDefault="00CC1722" w:rsidP="00CC1722"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>You must still be able to open and parse data elements and attributes from XML files (using TinyXML2 or similar).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(73).Range.InsertXML($paraXml)

# Paragraph 74: (5 points) GameConfig -> cyan
$paraXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6194D394" w14:textId="77777777" w:rsidR="00CF7E20" w:rsidRDefault="00CF7E20" w:rsidP="00616B75"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidR="00705E8E"><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t xml:space="preserve"> points) </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>GameConfig</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(74).Range.InsertXML($paraXml)

# Paragraph 75: You must support a GameConfig.txt (or GameConfig.xml) file... -> cyan
$paraXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2585AC49" w14:textId="77777777" w:rsidR="00CC1722" w:rsidRDefault="00CC1722" w:rsidP="00CC1722"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="3"/></w:numPr><w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t xml:space="preserve">You must support a “GameConfig.txt” (or GameConfig.xml) file whose values </w:t></w:r><w:r w:rsidR="0038420D"><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>the game can easily access to customize behavior.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(75).Range.InsertXML($paraXml)
